# Updated cryptos list on Tue Jan 30 05:35:41 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row on
# Sheet1, plus a few coins whose rank/row order shifted (Cosmos/Toncoin
# swap rows 29-30; EnergySwap/VeChain/Maker rotate rows 43-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.388.68'
$ws.Cells.Item(2, 5).Value = '  +2.82%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.308.16'
$ws.Cells.Item(3, 5).Value = '  +1.87%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''310.92'
$ws.Cells.Item(5, 5).Value = '  +1.50%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''102.02'
$ws.Cells.Item(6, 5).Value = '  +5.46%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.535'
$ws.Cells.Item(7, 5).Value = '  +1.47%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.00%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.527'
$ws.Cells.Item(9, 5).Value = '  +7.21%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''35.77'
$ws.Cells.Item(10, 5).Value = '  +2.13%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.0814'
$ws.Cells.Item(11, 5).Value = '  +2.93%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.78%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''6.99'
$ws.Cells.Item(13, 5).Value = '  +0.95%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.667.87'
$ws.Cells.Item(14, 5).Value = '  +1.91%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''14.99'
$ws.Cells.Item(15, 5).Value = '  +2.06%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.318.10'
$ws.Cells.Item(16, 5).Value = '  +2.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''0.809'
$ws.Cells.Item(17, 5).Value = '  +2.22%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '43.296.66'
$ws.Cells.Item(18, 5).Value = '  +2.89%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''12.31'
$ws.Cells.Item(19, 5).Value = '  -0.08%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0933'
$ws.Cells.Item(20, 5).Value = '  +3.07%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''6.18'
$ws.Cells.Item(21, 5).Value = '  +2.85%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''68.05'
$ws.Cells.Item(22, 5).Value = '  +0.40%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''241.40'
$ws.Cells.Item(23, 5).Value = '  +1.77%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''2.62'
$ws.Cells.Item(24, 5).Value = '  +1.94%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''2.00'
$ws.Cells.Item(25, 5).Value = '  +2.26%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.04%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''24.63'
$ws.Cells.Item(27, 5).Value = '  +4.70%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''36.98'
$ws.Cells.Item(28, 5).Value = '  -2.11%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).Value = '''2.19'
$ws.Cells.Item(29, 5).Value = '  +3.37%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(30, 4).Value = '''9.63'
$ws.Cells.Item(30, 5).Value = '  +0.50%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''168.07'
$ws.Cells.Item(31, 5).Value = '  +3.13%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''5.28'
$ws.Cells.Item(32, 5).Value = '  +0.75%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.09%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +5.93%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.0742'
$ws.Cells.Item(35, 5).Value = '  +0.65%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''17.60'
$ws.Cells.Item(36, 5).Value = '  +0.02%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''3.06'
$ws.Cells.Item(37, 5).Value = '  -3.40%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''1.88'
$ws.Cells.Item(38, 5).Value = '  +3.32%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +1.17%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +1.49%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''4.35'
$ws.Cells.Item(41, 5).Value = '  +7.40%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.51%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(43, 4).Value = '''0.0289'
$ws.Cells.Item(43, 5).Value = '  +2.94%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Maker'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(44, 4).Value = '1.972.83'
$ws.Cells.Item(44, 5).Value = '  +1.18%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).Value = '''19.29'
$ws.Cells.Item(45, 5).Value = '  +1.72%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''2.97'
$ws.Cells.Item(46, 5).Value = '  +1.93%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''9.90'
$ws.Cells.Item(47, 5).Value = '  +0.48%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''55.38'
$ws.Cells.Item(48, 5).Value = '  +2.67%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''2.90'
$ws.Cells.Item(49, 5).Value = '  +0.14%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +7.14%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '2.537.09'
$ws.Cells.Item(51, 5).Value = '  +1.88%  '
